# COM-interop edit script: recompute & update leve-profit figures
# across the Sagittarius_Profits workbook (scheduled-runner refresh).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3854.2778
$ws.Range("J17").Value = 3963.353
$ws.Range("L17").Value = 11890.059
$ws.Range("N17").Value = -12226.059
$ws.Range("H40").Value = 2140.7585
$ws.Range("I40").Value = 1974.4
$ws.Range("J40").Value = 2510.4443
$ws.Range("K40").Value = 1974.4
$ws.Range("L40").Value = 2510.4443
$ws.Range("M40").Value = -1799.4
$ws.Range("N40").Value = -2860.4443
$ws.Range("H42").Value = 1238.9166
$ws.Range("I42").Value = 1072.75
$ws.Range("J42").Value = 1571.25
$ws.Range("K42").Value = 3218.25
$ws.Range("L42").Value = 4713.75
$ws.Range("M42").Value = -2988.25
$ws.Range("N42").Value = -5173.75
$ws.Range("H43").Value = 1500
$ws.Range("I43").Value = 1500
$ws.Range("K43").Value = 1500
$ws.Range("M43").Value = -1431
$ws.Range("H100").Value = 4500
$ws.Range("I100").Value = 4000
$ws.Range("K100").Value = 4000
$ws.Range("M100").Value = -3459
$ws.Range("H112").Value = 3896.25
$ws.Range("J112").Value = 3995
$ws.Range("L112").Value = 11985
$ws.Range("N112").Value = -14201
$ws.Range("H127").Value = 5961
$ws.Range("I127").Value = 5961
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 17883
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -12923
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3350.6553
$ws.Range("I61").Value = 1367.5834
$ws.Range("K61").Value = 1367.5834
$ws.Range("M61").Value = -1155.5834
$ws.Range("H74").Value = 1816.6316
$ws.Range("J74").Value = 1544.4
$ws.Range("L74").Value = 1544.4
$ws.Range("N74").Value = -3292.4
$ws.Range("H77").Value = 1816.6316
$ws.Range("J77").Value = 1544.4
$ws.Range("L77").Value = 7722
$ws.Range("N77").Value = -16458
$ws.Range("H136").Value = 3350.6553
$ws.Range("I136").Value = 1367.5834
$ws.Range("K136").Value = 4102.7502
$ws.Range("M136").Value = -1552.7502

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1215.8334
$ws.Range("I86").Value = 1159
$ws.Range("K86").Value = 1159
$ws.Range("M86").Value = -36
$ws.Range("H89").Value = 1215.8334
$ws.Range("I89").Value = 1159
$ws.Range("K89").Value = 5795
$ws.Range("M89").Value = -179
$ws.Range("H105").Value = 3659.6
$ws.Range("I105").Value = 3659.6
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3659.6
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1912.6
$ws.Range("N105").ClearContents()
$ws.Range("H132").Value = 84499.5
$ws.Range("J132").Value = 84499.5
$ws.Range("L132").Value = 84499.5
$ws.Range("N132").Value = -94619.5
$ws.Range("H134").Value = 1237
$ws.Range("I134").Value = 1237
$ws.Range("K134").Value = 3711
$ws.Range("M134").Value = -1176

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2065.2307
$ws.Range("I31").Value = 1668.6666
$ws.Range("J31").Value = 2405.1428
$ws.Range("K31").Value = 1668.6666
$ws.Range("L31").Value = 2405.1428
$ws.Range("M31").Value = -1373.6666
$ws.Range("N31").Value = -2995.1428
$ws.Range("H34").Value = 2065.2307
$ws.Range("I34").Value = 1668.6666
$ws.Range("J34").Value = 2405.1428
$ws.Range("K34").Value = 1668.6666
$ws.Range("L34").Value = 2405.1428
$ws.Range("M34").Value = -1466.6666
$ws.Range("N34").Value = -2809.1428
$ws.Range("H99").Value = 3939.3845
$ws.Range("I99").Value = 3457.6
$ws.Range("J99").Value = 4240.5
$ws.Range("K99").Value = 3457.6
$ws.Range("L99").Value = 4240.5
$ws.Range("M99").Value = -1959.6
$ws.Range("N99").Value = -7236.5
$ws.Range("H126").Value = 3939.3845
$ws.Range("I126").Value = 3457.6
$ws.Range("J126").Value = 4240.5
$ws.Range("K126").Value = 10372.8
$ws.Range("L126").Value = 12721.5
$ws.Range("M126").Value = -7902.799999999999
$ws.Range("N126").Value = -17661.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 20028
$ws.Range("I44").Value = 20028
$ws.Range("K44").Value = 20028
$ws.Range("M44").Value = -19432
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H80").Value = 4733.1665
$ws.Range("J80").Value = 4950
$ws.Range("L80").Value = 4950
$ws.Range("N80").Value = -6946
$ws.Range("H83").Value = 4733.1665
$ws.Range("J83").Value = 4950
$ws.Range("L83").Value = 24750
$ws.Range("N83").Value = -34734
$ws.Range("H102").Value = 1013.9231
$ws.Range("I102").Value = 848.4167
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 848.4167
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = 773.5833
$ws.Range("N102").Value = -6244
$ws.Range("H132").Value = 1288.8
$ws.Range("I132").Value = 1123.8334
$ws.Range("J132").Value = 1536.25
$ws.Range("K132").Value = 3371.5002
$ws.Range("L132").Value = 4608.75
$ws.Range("M132").Value = -841.5001999999999
$ws.Range("N132").Value = -9668.75
$ws.Range("H134").Value = 274999.5
$ws.Range("J134").Value = 274999.5
$ws.Range("L134").Value = 824998.5
$ws.Range("N134").Value = -830068.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3858.625
$ws.Range("I40").Value = 2999
$ws.Range("K40").Value = 2999
$ws.Range("M40").Value = -2863
$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H122").Value = 5897.448
$ws.Range("I122").Value = 6012.3687
$ws.Range("K122").Value = 18037.1061
$ws.Range("M122").Value = -15587.1061
$ws.Range("H135").Value = 79329
$ws.Range("J135").Value = 79329
$ws.Range("L135").Value = 79329
$ws.Range("N135").Value = -89469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 17369.133
$ws.Range("J62").Value = 13065.429
$ws.Range("L62").Value = 13065.429
$ws.Range("N62").Value = -14313.429
$ws.Range("H65").Value = 17369.133
$ws.Range("J65").Value = 13065.429
$ws.Range("L65").Value = 65327.145
$ws.Range("N65").Value = -71567.145
$ws.Range("H96").Value = 1301.8
$ws.Range("I96").Value = 1049.5
$ws.Range("J96").Value = 1470
$ws.Range("K96").Value = 1049.5
$ws.Range("L96").Value = 1470
$ws.Range("M96").Value = 323.5
$ws.Range("N96").Value = -4216
$ws.Range("H126").Value = 4142.9443
$ws.Range("I126").Value = 2306.2727
$ws.Range("K126").Value = 6918.8181
$ws.Range("M126").Value = -4448.8181
$ws.Range("H132").Value = 3693.7273
$ws.Range("I132").Value = 4286.1763
$ws.Range("K132").Value = 12858.5289
$ws.Range("M132").Value = -10328.5289
